$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "women" image folder to "women1" for every Original/Scramble
# filename already on the sheet (in-place replace keeps shared-string slots
# stable instead of rewriting the whole table).
$ws.Cells.Replace("women/", "women1/")

# Add the newly completed experiment-2 reading (image pair 33).
$ws.Range("A34").Value = "women1/33-original.jpg"
$ws.Range("B34").Value = "women1/33-scramble.jpg"

# The filenames are one character longer now ("women" -> "women1"), so the
# best-fit columns need to grow to keep matching their content.
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(2).ColumnWidth = 23

# Turn on AutoFilter for the (now 34-row) Original column.
$ws.Range("A1:A34").AutoFilter(1)

# Excel's AutoFilter registers a hidden sheet-scoped _FilterDatabase name.
$fdb = $ws.Names.Add("_xlnm._FilterDatabase", "=Tabelle1!`$A`$1:`$A`$34")
$fdb.Visible = $false

# Move the active selection to B1.
$ws.Range("B1").Select()
